$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new shared strings in the exact order needed so the sharedStrings table
# ends up with: 18=Recherche, 19=Weekly-Summup-01 Meeting.,
# 20=Zusammenfassung Weekly-Summup-01 erstellt.,
# 21=Recherche bezüglich des Scrum Vorgehensmodells.
$ws.Range("C11").Value = "Recherche"
$ws.Range("D12").Value = "Weekly-Summup-01 Meeting."
$ws.Range("D13").Value = "Zusammenfassung Weekly-Summup-01 erstellt."
$ws.Range("D11").Value = "Recherche bezüglich des Scrum Vorgehensmodells."

# Reuse already-existing shared strings for the remaining new cells
$ws.Range("C12").Value = "Online-Meeting"
$ws.Range("C13").Value = "Planung"

# Dates (2024-10-21 -> serial 45586), formatted like the existing date column
$ws.Range("A11").Value = 45586
$ws.Range("A11").NumberFormat = "m/d/yy"
$ws.Range("A12").Value = 45586
$ws.Range("A12").NumberFormat = "m/d/yy"
$ws.Range("A13").Value = 45586
$ws.Range("A13").NumberFormat = "m/d/yy"

# Durations
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 1.25
$ws.Range("B13").Value = 1

# Match the style used by the rest of column B (left aligned)
$ws.Range("B11").HorizontalAlignment = -4131
$ws.Range("B12").HorizontalAlignment = -4131
$ws.Range("B13").HorizontalAlignment = -4131

# Update the selection to reflect the new "next empty row"
$ws.Range("A14").Select()
